# Generate Report for Handback
#
# The f5a43bc5-3649-4e54-ad56-505b15ae6eb9.md file has been handed back and is
# now in sync with en-US, so update its Status cells (which previously read
# "Ready for handoff") on the Overview, zh-cn and de-de sheets, and refresh the
# "Latest Handback DateTime" timestamps on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Overview sheet: row for f5a43bc5... is row 3 (zh-cn status in B3, de-de status in C3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# zh-cn detail sheet: row for f5a43bc5... is row 3 (Status in C3, Latest Handback DateTime in H3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("H3").Value = "2016-03-19 16:38:41"

# de-de detail sheet: row for f5a43bc5... is row 3 (Status in C3, Latest Handback DateTime in H3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("H3").Value = "2016-03-19 16:38:46"
